$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 118
$ws.Range("B118").Value = "Paulo Cesar Lucas Mendes"
$ws.Range("D118").Value = "Paulo Cesar Lucas Mendes"
$ws.Range("E118").Value = "KH7387741"
$ws.Range("F118").Value = "B"
$ws.Range("G118").Value = "R"
$ws.Range("H118").Value = "Mae - PQ5M-RCT / Pai GZF7-96M"

# Row 119
$ws.Range("A119").Value = "martinelly1@yahoo.com"
$ws.Range("B119").Value = "Martinelly Vieira Martins"
$ws.Range("C119").Value = "(31) 99559-7820"
$ws.Range("D119").Value = "Martinelly Vieira Martins"
$ws.Range("E119").Value = "UX8305420"
$ws.Range("F119").Value = "L1c3"
$ws.Range("H119").Value = "LYPX-6DF"

# Row 120
$ws.Range("A120").Value = "biancasantos.aps@gmail.com"
$ws.Range("B120").Value = "Bianca Mello"
$ws.Range("C120").Value = "(18) 996342066"
$ws.Range("D120").Value = "Bianca Mello"
$ws.Range("E120").Value = "GZ9117955"

# Row 121
$ws.Range("A121").Value = "renato.flister@gmail.com"
$ws.Range("B121").Value = "Renato Souza Lima"
$ws.Range("C121").Value = "31 99556-2131"
$ws.Range("D121").Value = "Renato Souza Lima"
$ws.Range("E121").Value = "CZ4140450"
$ws.Range("F121").Value = "L1c"
$ws.Range("G121").Value = "R1a"
$ws.Range("H121").Value = "Paterno = GR75-3DT / GR7P-2L7 - Materno = GQF1-4HB / GR75-ZJ6"
$ws.Range("L121").Value = "Salinas, Teofilo Otoni, Maranhao de Minas, Montes Claros, Topazio, Aracuiai, Pescador"
$ws.Range("M121").Value = "Pereira Lima, Schimidt, Roedel, Wolf, Both, Gomes Pereira, Souza Passos, Barbosa Lima, Francisco de Barros,Zeferino da Silva, Mata dos Santos, Franz"

# Row 122
$ws.Range("A122").Value = "biancasantos.aps@gmail.com"
$ws.Range("B122").Value = "Bianca Mello"
$ws.Range("C122").Value = "(18) 996342066"
$ws.Range("D122").Value = "Edeilza Barbosa dos Santos"
$ws.Range("E122").Value = "WE9043616"
$ws.Range("F122").Value = "C1"

# Update the selected/active cell to match the saved view state
$ws.Range("B109").Select()
